# Commit: fix bugs (order ds)
#
# A23 was an empty (but styled) cell; it now carries the text "sm".
# This also grows the shared-string table by one entry ("sm") and
# gives A23 its own cell-format record (same locked/hidden protection
# as the rest of the column) as Excel does when a previously blank,
# formatted cell receives a value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("A23")
$cell.Value = "sm"

# Preserve/assert the cell's protection formatting (locked, not hidden)
# that the rest of the column already uses.
$cell.Locked = $true
$cell.FormulaHidden = $false
